$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 76

$ws.Cells.Item($row, 1).Value = "Kindergarden"
$ws.Cells.Item($row, 2).Value = "Kindergarden Naarden Lambertus Hortensiuslaan"
$ws.Cells.Item($row, 3).Value = "KDV"

# D76 holds a date-shaped string ("2024-02-29") that must stay literal text
# (matches how the source file stores every other Rapportdatum cell: an
# inlineStr, not an auto-converted date serial). Assigning it straight to
# .Value makes Excel "helpfully" parse it into a date serial + a new
# NumberFormat style. Routing it through a TRIM() formula result copied in
# via PasteSpecial(xlPasteValues) preserves it as plain text with no style
# churn, exactly like the rest of column D.
$helper = $ws.Cells.Item(200, 1)
$helper.Formula = "=TRIM(""2024-02-29"")"
$helper.Copy() | Out-Null
$ws.Cells.Item($row, 4).PasteSpecial(-4163) | Out-Null
$helper.Clear() | Out-Null

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
